$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 520, shifting existing rows 520:593 down to 521:594
$ws.Rows.Item(520).Insert()

# Populate the newly inserted row 520 with the new weekly record
$ws.Cells.Item(520, 1).Value = 3
$ws.Cells.Item(520, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(520, 3).Value = "Coquimbo"
$ws.Cells.Item(520, 4).Value = 45131
$ws.Cells.Item(520, 5).Value = 5
$ws.Cells.Item(520, 6).Value = 100114013
$ws.Cells.Item(520, 7).Value = "Zanahoria"
$ws.Cells.Item(520, 8).Value = "Sin especificar"
$ws.Cells.Item(520, 9).Value = "Primera"
$ws.Cells.Item(520, 10).Value = 240
$ws.Cells.Item(520, 11).Value = 7000
$ws.Cells.Item(520, 12).Value = 7500
$ws.Cells.Item(520, 13).Value = 7250
$ws.Cells.Item(520, 14).Value = '$/saco 20 kilos'
$ws.Cells.Item(520, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(520, 16).Value = 362
$ws.Cells.Item(520, 17).Value = 20
$ws.Cells.Item(520, 18).Value = "Hortaliza"
